# Fix: rename the sheet back to "sem" and restore the row of values
# that was accidentally removed while fixing the spreadsheet size.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the worksheet (was reverted to the default "Sheet1" by mistake)
$ws.Name = "sem"

# 2) Restore the missing data row (row 2) under the header row
$ws.Range("A2").Value = "DOUBLES1"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 7
$ws.Range("E2").Value = "Ti-6Al-4V"
$ws.Range("F2").Value = 1.6
$ws.Range("G2").Value = 4.42
$ws.Range("I2").Value = 122.4135
$ws.Range("J2").Value = 105.1013

# Leave the selection on J2, matching the restored sheet's last cell
$ws.Range("J2").Select()
